$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.566.12'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '3.083.68'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.53%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.542'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.59%  '
$ws.Range("D9").Value = '3.079.76'
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("E10").Value = '  -1.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.83'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("E14").Value = '  -2.08%  '
$ws.Range("D15").Value = '3.597.91'
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").Value = '63.533.99'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").Value = '3.083.00'
$ws.Range("E19").Value = '  -1.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("E22").Value = '  -2.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  +3.63%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("E33").Value = '  +4.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("D35").Value = '0.0₃0846'
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.37'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '443.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.94%  '
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("E44").Value = '  -2.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("E46").Value = '  +3.26%  '
$ws.Range("D47").Value = '2.799.79'
$ws.Range("E47").Value = '  -3.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.68'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.63%  '
$ws.Range("E51").Value = '  +0.91%  '
